# Insert a new weekly price record as row 6, pushing the existing
# rows 6-8 down to rows 7-9 (dimension grows from A1:T8 to A1:T9).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(6).Insert()

$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C6").Value = 'Arica y Parinacota'
$ws.Range("D6").Value = 45071
$ws.Range("E6").Value = 15
$ws.Range("F6").Value = 'Fruta'
$ws.Range("G6").Value = 100107
$ws.Range("H6").Value = 'Otros'
$ws.Range("I6").Value = 100107001
$ws.Range("J6").Value = 'Caqui'
$ws.Range("K6").Value = 'Fuyu'
$ws.Range("L6").Value = 'Segunda'
$ws.Range("M6").Value = 110
$ws.Range("N6").Value = 23000
$ws.Range("O6").Value = 24000
$ws.Range("P6").Value = 23455
$ws.Range("Q6").Value = '$/caja 18 kilos granel'
$ws.Range("R6").Value = 'Región Metropolitana'
$ws.Range("S6").Value = 1303
$ws.Range("T6").Value = 18
